# This edit inserts one new data row into the "Hortaliza, Terminal La Palmera
# de La Serena - Ají" sheet. The new record is inserted right above the
# current row 520, pushing the existing rows 520-620 down to 521-621, and
# the sheet's used dimension grows from A1:R620 to A1:R621.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 520 (shifts rows 520:620 -> 521:621).
$ws.Rows.Item(520).Insert()

# Fill in the new row 520 with the new record's data.
$ws.Cells.Item(520, 1).Value  = 8
$ws.Cells.Item(520, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(520, 3).Value  = "Coquimbo"
$ws.Cells.Item(520, 4).Value  = 45258
$ws.Cells.Item(520, 5).Value  = 4
$ws.Cells.Item(520, 6).Value  = 100112021
$ws.Cells.Item(520, 7).Value  = "Ají"
$ws.Cells.Item(520, 8).Value  = "Inferno"
$ws.Cells.Item(520, 9).Value  = "Primera"
$ws.Cells.Item(520, 10).Value = 440
$ws.Cells.Item(520, 11).Value = 29000
$ws.Cells.Item(520, 12).Value = 30000
$ws.Cells.Item(520, 13).Value = 29500
$ws.Cells.Item(520, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(520, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(520, 16).Value = 2950
$ws.Cells.Item(520, 17).Value = 10
$ws.Cells.Item(520, 18).Value = "Hortaliza"
